$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append duration annotations to three task labels ---
$ws.Range("B56").Value = "choixConnect.html (20mn)"
$ws.Range("B64").Value = "logUser.php (1h)"
$ws.Range("B75").Value = "traitLogUser.php (30mn)"

# B56 picks up the bold "section sub-item" look shared by its siblings
# (B29 / B43 already use this bold style) instead of its previous plain style.
$ws.Range("B56").Font.Bold = $true

# --- Mark the matching "done" column (C) with "ok" for these rows ---
$doneRows = @(59, 62, 67, 68, 74, 81, 82, 84, 85)
foreach ($r in $doneRows) {
    $ws.Range("C$r").Value = "ok"
}

# --- Restore the view/selection state recorded in the saved workbook ---
$ws.Application.ActiveWindow.ScrollRow = 55
$ws.Range("J71").Select() | Out-Null
